# Add a new bullet item right after the paragraph that ends with
# "...uzlikt atpūtas laiku starp vingrinājumiem." describing the new
# registration / login functionality.

$d = $word.ActiveDocument

# Locate the paragraph to insert after by its distinctive trailing text.
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("starp vingrinājumiem.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorPara = $anchorRange.Paragraphs(1)

# Insert a brand-new paragraph right after it; it inherits the same
# pPr (numbered-list formatting) as the anchor paragraph, matching the
# existing bullet items in the document.
$newRange = $anchorPara.Range.InsertParagraphAfter()

# Re-fetch the freshly created (still empty) paragraph and fill it in.
$newPara = $anchorPara.Next()
$newPara.Range.InsertAfter("Lietotājam būs iespēja piereģistrēties un izveidot kontu ar kuru varēs ieejiet aplikācija(Reģistrācija būs vajadzīgs: Full name, epasts, username, parole. Priekš Login būs vajadzīgs username un parole)")
